# Updated cryptos list (prices + 1h volume deltas); rows 31/32 also swap Coin/Link.
# D-column values that look like plain numbers ("0.9997", "241.86", ...) are written
# with a leading apostrophe so Excel stores them as literal text (matching the source
# data, which keeps multi-dot price strings like "29.307.95" as text too) instead of
# silently parsing them into numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.307.95'
$ws.Range("E2").Value = '  +0.32%  '

# Row 3
$ws.Range("D3").Value = '1.874.83'
$ws.Range("E3").Value = '  +0.54%  '

# Row 4
$ws.Range("D4").Value = '''0.9993'
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").Value = '''0.7132'
$ws.Range("E5").Value = '  -0.27%  '

# Row 6
$ws.Range("D6").Value = '''241.86'
$ws.Range("E6").Value = '  +0.55%  '

# Row 7
$ws.Range("D7").Value = '''0.9997'
$ws.Range("E7").Value = '  -0.27%  '

# Row 8
$ws.Range("D8").Value = '''0.3107'
$ws.Range("E8").Value = '  +1.20%  '

# Row 9
$ws.Range("D9").Value = '''0.07725'
$ws.Range("E9").Value = '  -0.31%  '

# Row 10
$ws.Range("D10").Value = '''25.07'
$ws.Range("E10").Value = '  +0.57%  '

# Row 11
$ws.Range("D11").Value = '''0.08385'
$ws.Range("E11").Value = '  +1.66%  '

# Row 12
$ws.Range("D12").Value = '1.884.26'
$ws.Range("E12").Value = '  +0.90%  '

# Row 13
$ws.Range("D13").Value = '''5.217'
$ws.Range("E13").Value = '  +0.15%  '

# Row 14
$ws.Range("D14").Value = '''0.7114'
$ws.Range("E14").Value = '  -0.62%  '

# Row 15
$ws.Range("D15").Value = '''91.35'
$ws.Range("E15").Value = '  +1.21%  '

# Row 16
$ws.Range("D16").Value = '29.302.41'
$ws.Range("E16").Value = '  +0.29%  '

# Row 17
$ws.Range("D17").Value = '''0.000008272'
$ws.Range("E17").Value = '  +6.38%  '

# Row 18
$ws.Range("D18").Value = '''5.977'
$ws.Range("E18").Value = '  +2.50%  '

# Row 19
$ws.Range("D19").Value = '''242.69'
$ws.Range("E19").Value = '  -0.05%  '

# Row 20
$ws.Range("D20").Value = '2.129.14'
$ws.Range("E20").Value = '  +0.14%  '

# Row 21
$ws.Range("E21").Value = '  +0.78%  '

# Row 22
$ws.Range("D22").Value = '''0.9992'
$ws.Range("E22").Value = '  -0.29%  '

# Row 23
$ws.Range("D23").Value = '''7.825'
$ws.Range("E23").Value = '  -1.51%  '

# Row 24
$ws.Range("D24").Value = '''0.9997'
$ws.Range("E24").Value = '  -0.26%  '

# Row 25
$ws.Range("E25").Value = '  +2.05%  '

# Row 26
$ws.Range("E26").Value = '  +0.77%  '

# Row 27
$ws.Range("D27").Value = '''9.021'
$ws.Range("E27").Value = '  +1.37%  '

# Row 28
$ws.Range("D28").Value = '''18.53'
$ws.Range("E28").Value = '  +2.14%  '

# Row 29
$ws.Range("D29").Value = '''1.504'
$ws.Range("E29").Value = '  +0.71%  '

# Row 30
$ws.Range("D30").Value = '''4.421'
$ws.Range("E30").Value = '  +1.84%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''4.326'
$ws.Range("E31").Value = '  +6.12%  '

# Row 32
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D32").Value = '''1.289'
$ws.Range("E32").Value = '  -1.07%  '

# Row 33
$ws.Range("D33").Value = '''0.05246'
$ws.Range("E33").Value = '  +1.14%  '

# Row 34
$ws.Range("D34").Value = '''1.927'
$ws.Range("E34").Value = '  +0.75%  '

# Row 35
$ws.Range("D35").Value = '''0.7481'
$ws.Range("E35").Value = '  +2.87%  '

# Row 36
$ws.Range("D36").Value = '''1.173'
$ws.Range("E36").Value = '  -0.15%  '

# Row 37
$ws.Range("D37").Value = '''2.682'

# Row 38
$ws.Range("D38").Value = '''0.01858'
$ws.Range("E38").Value = '  +0.73%  '

# Row 39
$ws.Range("E39").Value = '  +0.96%  '

# Row 40
$ws.Range("D40").Value = '1.156.19'
$ws.Range("E40").Value = '  -0.81%  '

# Row 41
$ws.Range("D41").Value = '''6.368'

# Row 42
$ws.Range("E42").Value = '  +1.32%  '

# Row 43
$ws.Range("D43").Value = '''0.8853'
$ws.Range("E43").Value = '  -1.92%  '

# Row 44
$ws.Range("D44").Value = '''105.64'
$ws.Range("E44").Value = '  +3.85%  '

# Row 45
$ws.Range("D45").Value = '''0.9994'
$ws.Range("E45").Value = '  -0.34%  '

# Row 46
$ws.Range("D46").Value = '2.025.83'
$ws.Range("E46").Value = '  +0.36%  '

# Row 47
$ws.Range("D47").Value = '''1.804'
$ws.Range("E47").Value = '  +2.54%  '

# Row 48
$ws.Range("D48").Value = '''0.5189'
$ws.Range("E48").Value = '  -1.81%  '

# Row 49
$ws.Range("E49").Value = '  +4.30%  '

# Row 50
$ws.Range("D50").Value = '''9.383'
$ws.Range("E50").Value = '  +1.44%  '

# Row 51
$ws.Range("D51").Value = '''0.4301'
$ws.Range("E51").Value = '  +1.68%  '
